$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.868.03'
$ws.Range('D3').Value = '1.630.11'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5069'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2577'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06334'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.45'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.640.10'
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.249'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').Value = '1.854.16'
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5507'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('D18').Value = '25.884.80'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.403'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.880'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.024'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('E25').Value = '  +2.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('E27').Value = '  +3.50%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.758'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('E34').Value = '  +1.18%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8961'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5523'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.537'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.39%  '
$ws.Range('D39').Value = '1.118.94'
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01555'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.71%  '
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.587'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7977'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.19'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.97%  '
$ws.Range('E45').Value = '  -2.93%  '
$ws.Range('D46').Value = '1.765.24'
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4443'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.17%  '
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05135'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.582'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.82%  '
